$d = $word.ActiveDocument

# Update the header date line
$d.Content.Find.Execute("2024-01-08 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-09 Tuesday", 2) | Out-Null

# Update each multiplication-fact cell in the practice table.
# wdReplaceOne (1), scoped to a single cell's Range, is used throughout
# because some old values (e.g. "57x23=1311") are duplicated across
# cells and each occurrence must map to a different new value.
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Find.Execute("79×35=2765", $true, $false, $false, $false, $false, $true, 1, $false, "86×71=6106", 1) | Out-Null
$t.Cell(1, 2).Range.Find.Execute("28×83=2324", $true, $false, $false, $false, $false, $true, 1, $false, "94×53=4982", 1) | Out-Null
$t.Cell(1, 3).Range.Find.Execute("89×71=6319", $true, $false, $false, $false, $false, $true, 1, $false, "55×90=4950", 1) | Out-Null
$t.Cell(1, 4).Range.Find.Execute("59×85=5015", $true, $false, $false, $false, $false, $true, 1, $false, "98×50=4900", 1) | Out-Null
$t.Cell(1, 5).Range.Find.Execute("61×49=2989", $true, $false, $false, $false, $false, $true, 1, $false, "85×96=8160", 1) | Out-Null
$t.Cell(5, 1).Range.Find.Execute("36×65=2340", $true, $false, $false, $false, $false, $true, 1, $false, "33×63=2079", 1) | Out-Null
$t.Cell(5, 2).Range.Find.Execute("35×65=2275", $true, $false, $false, $false, $false, $true, 1, $false, "32×35=1120", 1) | Out-Null
$t.Cell(5, 3).Range.Find.Execute("14×65=910", $true, $false, $false, $false, $false, $true, 1, $false, "43×85=3655", 1) | Out-Null
$t.Cell(5, 4).Range.Find.Execute("34×78=2652", $true, $false, $false, $false, $false, $true, 1, $false, "22×81=1782", 1) | Out-Null
$t.Cell(5, 5).Range.Find.Execute("72×23=1656", $true, $false, $false, $false, $false, $true, 1, $false, "31×72=2232", 1) | Out-Null
$t.Cell(10, 1).Range.Find.Execute("11×67=737", $true, $false, $false, $false, $false, $true, 1, $false, "39×67=2613", 1) | Out-Null
$t.Cell(10, 2).Range.Find.Execute("57×23=1311", $true, $false, $false, $false, $false, $true, 1, $false, "71×37=2627", 1) | Out-Null
$t.Cell(10, 3).Range.Find.Execute("28×97=2716", $true, $false, $false, $false, $false, $true, 1, $false, "12×71=852", 1) | Out-Null
$t.Cell(10, 4).Range.Find.Execute("36×57=2052", $true, $false, $false, $false, $false, $true, 1, $false, "43×87=3741", 1) | Out-Null
$t.Cell(10, 5).Range.Find.Execute("54×58=3132", $true, $false, $false, $false, $false, $true, 1, $false, "94×37=3478", 1) | Out-Null
$t.Cell(15, 1).Range.Find.Execute("57×23=1311", $true, $false, $false, $false, $false, $true, 1, $false, "46×49=2254", 1) | Out-Null
$t.Cell(15, 2).Range.Find.Execute("47×81=3807", $true, $false, $false, $false, $false, $true, 1, $false, "40×74=2960", 1) | Out-Null
$t.Cell(15, 3).Range.Find.Execute("32×40=1280", $true, $false, $false, $false, $false, $true, 1, $false, "57×28=1596", 1) | Out-Null
$t.Cell(15, 4).Range.Find.Execute("53×29=1537", $true, $false, $false, $false, $false, $true, 1, $false, "41×23=943", 1) | Out-Null
$t.Cell(15, 5).Range.Find.Execute("25×33=825", $true, $false, $false, $false, $false, $true, 1, $false, "58×69=4002", 1) | Out-Null
$t.Cell(20, 1).Range.Find.Execute("91×67=6097", $true, $false, $false, $false, $false, $true, 1, $false, "86×19=1634", 1) | Out-Null
$t.Cell(20, 2).Range.Find.Execute("84×76=6384", $true, $false, $false, $false, $false, $true, 1, $false, "62×45=2790", 1) | Out-Null
$t.Cell(20, 3).Range.Find.Execute("55×12=660", $true, $false, $false, $false, $false, $true, 1, $false, "43×66=2838", 1) | Out-Null
$t.Cell(20, 4).Range.Find.Execute("91×81=7371", $true, $false, $false, $false, $false, $true, 1, $false, "46×34=1564", 1) | Out-Null
$t.Cell(20, 5).Range.Find.Execute("71×42=2982", $true, $false, $false, $false, $false, $true, 1, $false, "65×45=2925", 1) | Out-Null
